# Apply updated cryptocurrency price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '37.621.79'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.32%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.084.65'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.37%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '233.75'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.07%  '
$ws.Range("E6").Value = '  +2.09%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '58.04'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.34%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.392'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.72%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0781'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.52%  '
$ws.Range("E11").Value = '  +3.23%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.17'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.26%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.390.49'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.33%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.774'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.11%  '
$ws.Range("E16").Value = '  +1.18%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.077.88'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.30%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '37.587.78'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.21%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.06'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.90%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '70.84'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.42%  '
$ws.Range("E21").Value = '  +0.00%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '229.26'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.30%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.03%  '
$ws.Range("E24").Value = '  -0.60%  '
$ws.Range("E25").Value = '  -2.56%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.72'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +7.59%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '170.79'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.25%  '
$ws.Range("E28").Value = '  -3.87%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.50'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.20%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.39'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.35%  '
$ws.Range("E31").Value = '  +1.01%  '
$ws.Range("E32").Value = '  +0.08%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0639'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.41%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.67'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.77%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.50'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.24%  '
$ws.Range("E36").Value = '  -0.79%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.34'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.50%  '
$ws.Range("E38").Value = '  +0.07%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.38'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.27%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0234'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +8.77%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '100.96'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.02%  '
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.20'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +4.44%  '
$ws.Range("B43").Value = 'Cronos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0959'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.05%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.91'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.92%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '16.94'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.84%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.465.30'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.10%  '
$ws.Range("E47").Value = '  -0.58%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.02'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -5.43%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.25'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.05%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.95'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.85%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.273.56'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.27%  '
